# Updates the crypto price table on Sheet1 with refreshed Price (column D)
# and Volume(1h) (column E) figures, matching the GitHub Actions data
# refresh. Rows 19-20 and 22-23 also swap their Coin/Link/Price/Volume
# content (ranking re-order between Chainlink/WrappedEther and
# Polygon/Uniswap).
#
# Numeric-looking Price values (e.g. "16.65") are written with a leading
# apostrophe so Excel stores them as text (matching the source data,
# which keeps these as literal strings, not numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.337.32'
$ws.Range("E2").Value = '  -0.81%  '
$ws.Range("D3").Value = '3.076.06'
$ws.Range("E3").Value = '  -1.46%  '
$ws.Range("D5").Value = '''575.34'
$ws.Range("E5").Value = '  -0.49%  '
$ws.Range("D6").Value = '''171.06'
$ws.Range("E6").Value = '  -0.52%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = '3.073.22'
$ws.Range("E8").Value = '  -1.42%  '
$ws.Range("E9").Value = '  -2.10%  '
$ws.Range("D10").Value = '''6.24'
$ws.Range("E10").Value = '  -2.28%  '
$ws.Range("E11").Value = '  -2.35%  '
$ws.Range("E12").Value = '  -2.52%  '
$ws.Range("E13").Value = '  -3.99%  '
$ws.Range("D14").Value = '''35.78'
$ws.Range("E14").Value = '  -3.63%  '
$ws.Range("E15").Value = '  -1.41%  '
$ws.Range("D16").Value = '3.587.94'
$ws.Range("E16").Value = '  -1.38%  '
$ws.Range("D17").Value = '66.316.30'
$ws.Range("E17").Value = '  -0.82%  '
$ws.Range("E18").Value = '  -2.54%  '
$ws.Range("B19").Value = 'Chainlink'
$ws.Range("C19").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D19").Value = '''16.65'
$ws.Range("E19").Value = '  +2.38%  '
$ws.Range("B20").Value = 'WrappedEther'
$ws.Range("C20").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D20").Value = '3.078.54'
$ws.Range("E20").Value = '  -1.32%  '
$ws.Range("D21").Value = '''487.05'
$ws.Range("E21").Value = '  +2.62%  '
$ws.Range("B22").Value = 'Polygon'
$ws.Range("C22").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D22").Value = '''0.686'
$ws.Range("E22").Value = '  -3.25%  '
$ws.Range("B23").Value = 'Uniswap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D23").Value = '''7.68'
$ws.Range("E23").Value = '  -2.46%  '
$ws.Range("D24").Value = '''82.42'
$ws.Range("E24").Value = '  -1.56%  '
$ws.Range("D25").Value = '''12.63'
$ws.Range("E25").Value = '  -4.85%  '
$ws.Range("D26").Value = '''2.21'
$ws.Range("E26").Value = '  -2.97%  '
$ws.Range("D27").Value = '''10.11'
$ws.Range("E27").Value = '  -1.75%  '
$ws.Range("D28").Value = '''0.999'
$ws.Range("E28").Value = '  -0.10%  '
$ws.Range("D29").Value = '''7.85'
$ws.Range("E29").Value = '  -0.64%  '
$ws.Range("E30").Value = '  -5.19%  '
$ws.Range("E31").Value = '  -3.19%  '
$ws.Range("D32").Value = '''27.73'
$ws.Range("E32").Value = '  -3.13%  '
$ws.Range("E33").Value = '  -3.70%  '
$ws.Range("D34").Value = '0.0₃0910'
$ws.Range("E34").Value = '  -4.15%  '
$ws.Range("D35").Value = '''1.00'
$ws.Range("E35").Value = '  +0.06%  '
$ws.Range("D36").Value = '''47.93'
$ws.Range("E36").Value = '  +1.92%  '
$ws.Range("E37").Value = '  -4.66%  '
$ws.Range("E38").Value = '  -3.36%  '
$ws.Range("D39").Value = '''0.123'
$ws.Range("E39").Value = '  -0.88%  '
$ws.Range("E40").Value = '  -3.20%  '
$ws.Range("E41").Value = '  -4.65%  '
$ws.Range("E42").Value = '  -4.37%  '
$ws.Range("D43").Value = '2.772.61'
$ws.Range("E43").Value = '  -1.65%  '
$ws.Range("E45").Value = '  -2.66%  '
$ws.Range("D46").Value = '''134.76'
$ws.Range("E46").Value = '  -0.67%  '
$ws.Range("D47").Value = '''365.25'
$ws.Range("E47").Value = '  -4.59%  '
$ws.Range("D49").Value = '''24.22'
$ws.Range("E49").Value = '  -2.90%  '
$ws.Range("E50").Value = '  -2.34%  '
$ws.Range("E51").Value = '  -2.09%  '
